$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix DB logic, ratings, and handlers
$ws.Range("C2").Value = 4.78
$ws.Range("D2").Value = 27

$ws.Range("C3").Value = 4.4
$ws.Range("D3").Value = 15
